$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells keep their original text formatting
# (price strings must not be auto-converted to numbers by Excel)

$ws.Range("D2").Value = '65.962.14'
$ws.Range("E2").Value = '  -1.14%  '

$ws.Range("D3").Value = '3.449.28'
$ws.Range("E3").Value = '  -0.21%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.73'
$ws.Range("E5").Value = '  +0.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.95'
$ws.Range("E6").Value = '  -0.92%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.600'
$ws.Range("E8").Value = '  -0.21%  '

$ws.Range("D9").Value = '3.447.46'
$ws.Range("E9").Value = '  -0.26%  '

$ws.Range("E10").Value = '  -2.11%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.93'
$ws.Range("E11").Value = '  +1.05%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.409'
$ws.Range("E12").Value = '  -2.87%  '

$ws.Range("D13").Value = '4.042.42'
$ws.Range("E13").Value = '  -0.23%  '

$ws.Range("E14").Value = '  +1.64%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.83'
$ws.Range("E15").Value = '  -7.20%  '

$ws.Range("D16").Value = '65.930.96'
$ws.Range("E16").Value = '  -1.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000171'
$ws.Range("E17").Value = '  -0.89%  '

$ws.Range("D18").Value = '3.432.89'
$ws.Range("E18").Value = '  -0.64%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.94'
$ws.Range("E19").Value = '  -1.50%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.79'
$ws.Range("E20").Value = '  -0.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '369.08'
$ws.Range("E21").Value = '  -1.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.61'
$ws.Range("E22").Value = '  -1.86%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.26'
$ws.Range("E23").Value = '  +1.67%  '

$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("E25").Value = '  +0.61%  '

$ws.Range("E26").Value = '  +2.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.70'
$ws.Range("E27").Value = '  -1.93%  '

$ws.Range("E28").Value = '  +2.80%  '

$ws.Range("E29").Value = '  +0.39%  '

$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.77'
$ws.Range("E30").Value = '  -2.18%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '23.68'
$ws.Range("E31").Value = '  -1.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.98'
$ws.Range("E32").Value = '  -1.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.00%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.29'
$ws.Range("E34").Value = '  -4.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.98'
$ws.Range("E35").Value = '  -2.15%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.53'
$ws.Range("E36").Value = '  +0.74%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.95'
$ws.Range("E37").Value = '  +0.35%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.881'
$ws.Range("E38").Value = '  +0.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '28.75'
$ws.Range("E39").Value = '  +4.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.78'
$ws.Range("E40").Value = '  -0.92%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.63'
$ws.Range("E41").Value = '  -0.48%  '

$ws.Range("D42").Value = '2.759.31'
$ws.Range("E42").Value = '  +2.47%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.45'
$ws.Range("E43").Value = '  +0.17%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.44'
$ws.Range("E44").Value = '  -1.44%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0678'
$ws.Range("E45").Value = '  -2.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.25'
$ws.Range("E46").Value = '  -0.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.66'
$ws.Range("E47").Value = '  -2.37%  '

$ws.Range("E48").Value = '  -1.58%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '326.42'
$ws.Range("E49").Value = '  +1.63%  '

$ws.Range("E50").Value = '  -0.47%  '

$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '31.46'
$ws.Range("E51").Value = '  +0.21%  '
